$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "strongyl"
$ws.Range("A2").Select()
